$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as text
# (matching the source data which is all shared-string/text), then reset
# the cell style back to Normal so no stray number-format style sticks
# around on the cell (keeps cell's "s" attribute as the default style).
function Set-TextCell($row, $col, $val) {
    $target = $ws.Cells.Item($row, $col)
    $target.NumberFormat = "@"
    $target.Value = $val
    $target.Style = "Normal"
}

$data = @(
    @('University Extra','20','126305','07-04-2025','10:00 AM TO 01:00 PM','1','1'),
    @('University Extra','20','303001','08-04-2025','02:00 PM TO 05:00 PM','1','1'),
    @('University Extra','20','303002','09-04-2025','02:00 PM TO 05:00 PM','1','1'),
    @('University Extra','20','303015','17-04-2025','02:00 PM TO 05:00 PM','1','1'),
    @('University Extra','20','310021','07-04-2025','10:00 AM TO 01:00 PM','1','1'),
    @('University Extra','20','311511','07-04-2025','10:00 AM TO 01:00 PM','1','1'),
    @('University Extra','20','321015','07-04-2025','02:00 PM TO 05:00 PM','1','1'),
    @('University Extra','20','321125','07-04-2025','02:00 PM TO 05:00 PM','1','1'),
    @('University Extra','20','321165','08-04-2025','10:00 AM TO 01:00 PM','1','1'),
    @('University Extra','20','321175','08-04-2025','10:00 AM TO 01:00 PM','1','1'),
    @('University Extra','20','321235','07-04-2025','02:00 PM TO 05:00 PM','1','1')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 1
        Set-TextCell $row $col $rowData[$j]
    }
}

$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
